# Re-sort the facilities table (Table1, A1:H26) by the "#fid" column (A)
# ascending instead of the previous sort by "category" (H), and update the
# active selection to reflect where the user clicked afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range by column A (#fid), keeping the header row in place.
$dataRange = $ws.Range("A1:H26")
$sortKey   = $ws.Range("A1:A26")
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)

# Reflect the new cell selection left after the re-sort / review.
$ws.Range("P5").Select() | Out-Null
